$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells are treated as text so values such as "1.00" or
# "233.28" are not re-interpreted as numbers/dates by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.375.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.977.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -11.08%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.38%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.55"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.83"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.371"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0985"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.268.71"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.95"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.96"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.749"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -7.27%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.973.04"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.357.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0803"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.52"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -11.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.125"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.36"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.30%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.41"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.96%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.456.54"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0901"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0201"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.89%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -10.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.92"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.80"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.989"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.74"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.161.45"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.03%  "
